$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vendas")

# O285: add a new "Pagar Sapato" value of 20 (was blank)
$ws.Range("O285").Value = 20

# Row 286: Caterpilhar Cano Curto sale
$ws.Range("A286").Value = 45784
$ws.Range("B286").Value = "Loja2"
$ws.Range("C286").Value = 2188
$ws.Range("D286").Value = 38
$ws.Range("E286").Value = "Netony"
$ws.Range("F286").Value = "Caterpilhar Cano Curto"
$ws.Range("G286").Value = 100
$ws.Range("H286").Value = 199
$ws.Range("I286").Value = 1
$ws.Range("J286").Value = 180

# Row 287: Sapatilha sale
$ws.Range("A287").Value = 45784
$ws.Range("B287").Value = "Loja2"
$ws.Range("C287").Value = 1046
$ws.Range("D287").Value = 39
$ws.Range("E287").Value = "Rossanfort"
$ws.Range("F287").Value = "Sapatilha"
$ws.Range("G287").Value = 50
$ws.Range("H287").Value = 110
$ws.Range("I287").Value = 1
$ws.Range("J287").Value = 100

# Row 288: Sapatenis sale
$ws.Range("A288").Value = 45784
$ws.Range("B288").Value = "Loja2"
$ws.Range("C288").Value = 5500
$ws.Range("D288").Value = 39
$ws.Range("E288").Value = "Netony"
$ws.Range("F288").Value = "Sapatenis"
$ws.Range("G288").Value = 112
$ws.Range("H288").Value = 190
$ws.Range("I288").Value = 1
$ws.Range("J288").Value = 190

# Row 289: Tenis sale
$ws.Range("A289").Value = 45784
$ws.Range("B289").Value = "Loja2"
$ws.Range("C289").Value = 400
$ws.Range("D289").Value = 40
$ws.Range("E289").Value = "Ranster"
$ws.Range("F289").Value = "Tenis"
$ws.Range("G289").Value = 126
$ws.Range("H289").Value = 210
$ws.Range("I289").Value = 1
$ws.Range("J289").Value = 210

# Row 291: remove the leftover "teste" scratch note (W291/X291)
$ws.Range("W291:X291").Clear()

# Highlight the growing shortfall total in bold red, like the other blocks
$ws.Range("W304").Font.Bold = $true
$ws.Range("W304").Font.Color = 255

# Update the view to where the new rows were entered
$ws.Activate()
$ws.Range("X302").Select()
